$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of matchup data to append below the existing data (rows 927-942)
$data = @(
    @(3, 13, 5, 7),
    @(2, 14, 3, 6),
    @(4, 6, 6, 14),
    @(4, 16, 5, 4),
    @(4, 17, 5, 3),
    @(4, 7, 5, 13),
    @(5, 5, 4, 15),
    @(5, 15, 7, 5),
    @(4, 18, 3, 2),
    @(3, 5, 4, 15),
    @(4, 6, 3, 14),
    @(6, 16, 4, 4),
    @(5, 18, 6, 2),
    @(5, 8, 4, 12),
    @(7, 16, 4, 4),
    @(4, 13, 2, 7)
)

$startRow = 927
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

# Update the sheet view to match the new data extent
$ws.Application.ActiveWindow.ScrollRow = 932
$ws.Range("A943").Select()
